# Scheduled runner update: refresh cached market-board price / profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) on the various
# crafting-job leve sheets. Values below were pulled from the latest
# market data snapshot; some rows also gain/lose a LeveProfitNQ (M) or
# LeveProfitHQ (N) cell depending on whether that recipe currently has
# an NQ/HQ price quote available.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 312.625
$ws.Range("I2").Value = 312.625
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 312.625
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -199.625
$ws.Range("N2").ClearContents()

$ws.Range("H6").Value = 1401.625
$ws.Range("J6").Value = 5499.5
$ws.Range("L6").Value = 16498.5
$ws.Range("N6").Value = -16722.5

$ws.Range("H12").Value = 1573.0435
$ws.Range("I12").Value = 1009.25
$ws.Range("K12").Value = 1009.25
$ws.Range("M12").Value = -839.25

$ws.Range("H15").Value = 479.74075
$ws.Range("I15").Value = 479.74075
$ws.Range("K15").Value = 1439.22225
$ws.Range("M15").Value = -1270.22225

$ws.Range("H62").Value = 7614.3335
$ws.Range("I62").Value = 3999
$ws.Range("K62").Value = 3999
$ws.Range("M62").Value = -3375

$ws.Range("H65").Value = 7614.3335
$ws.Range("I65").Value = 3999
$ws.Range("K65").Value = 19995
$ws.Range("M65").Value = -16875

$ws.Range("H132").Value = 4910.8
$ws.Range("I132").Value = 2622.2
$ws.Range("K132").Value = 7866.599999999999
$ws.Range("M132").Value = -5336.599999999999

$ws.Range("H135").Value = 807.1
$ws.Range("I135").Value = 807.1
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7263.900000000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4728.900000000001
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 2049.2144
$ws.Range("I137").Value = 1711.125
$ws.Range("K137").Value = 5133.375
$ws.Range("M137").Value = -2583.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 54999
$ws.Range("J23").Value = 54999
$ws.Range("L23").Value = 54999
$ws.Range("N23").Value = -55517

$ws.Range("H32").Value = 4104.6763
$ws.Range("I32").Value = 4104.6763
$ws.Range("K32").Value = 4104.6763
$ws.Range("M32").Value = -3817.6763

$ws.Range("H45").Value = 4722.375
$ws.Range("I45").Value = 3441.25
$ws.Range("K45").Value = 3441.25
$ws.Range("M45").Value = -3064.25

$ws.Range("H61").Value = 3198.111
$ws.Range("I61").Value = 3198.111
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3198.111
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2986.111
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 791.6667
$ws.Range("I74").Value = 790
$ws.Range("J74").Value = 800
$ws.Range("K74").Value = 790
$ws.Range("L74").Value = 800
$ws.Range("M74").Value = 84
$ws.Range("N74").Value = -2548

$ws.Range("H77").Value = 791.6667
$ws.Range("I77").Value = 790
$ws.Range("J77").Value = 800
$ws.Range("K77").Value = 3950
$ws.Range("L77").Value = 4000
$ws.Range("M77").Value = 418
$ws.Range("N77").Value = -12736

$ws.Range("H97").Value = 8722.857
$ws.Range("I97").Value = 195.5
$ws.Range("K97").Value = 195.5
$ws.Range("M97").Value = 300.5

$ws.Range("H102").Value = 2346.5
$ws.Range("I102").Value = 769.875
$ws.Range("K102").Value = 769.875
$ws.Range("M102").Value = 852.125

$ws.Range("H132").Value = 1330.7
$ws.Range("I132").Value = 1330.7
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3992.1
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1462.1
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 3198.111
$ws.Range("I136").Value = 3198.111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9594.332999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7044.332999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 410.77777
$ws.Range("I80").Value = 365
$ws.Range("K80").Value = 365
$ws.Range("M80").Value = 633

$ws.Range("H83").Value = 410.77777
$ws.Range("I83").Value = 365
$ws.Range("K83").Value = 1825
$ws.Range("M83").Value = 3167

$ws.Range("H86").Value = 2143.3635
$ws.Range("I86").Value = 1968.2858
$ws.Range("J86").Value = 2449.75
$ws.Range("K86").Value = 1968.2858
$ws.Range("L86").Value = 2449.75
$ws.Range("M86").Value = -845.2858000000001
$ws.Range("N86").Value = -4695.75

$ws.Range("H89").Value = 2143.3635
$ws.Range("I89").Value = 1968.2858
$ws.Range("J89").Value = 2449.75
$ws.Range("K89").Value = 9841.429
$ws.Range("L89").Value = 12248.75
$ws.Range("M89").Value = -4225.429
$ws.Range("N89").Value = -23480.75

$ws.Range("H105").Value = 4224.6665
$ws.Range("I105").Value = 4224.6665
$ws.Range("K105").Value = 4224.6665
$ws.Range("M105").Value = -2477.6665

$ws.Range("H134").Value = 3782.8333
$ws.Range("I134").Value = 3043.7778
$ws.Range("K134").Value = 9131.3334
$ws.Range("M134").Value = -6596.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1123.5652
$ws.Range("I105").Value = 757.2778
$ws.Range("K105").Value = 757.2778
$ws.Range("M105").Value = 989.7222

$ws.Range("H132").Value = 5344.3887
$ws.Range("I132").Value = 5344.3887
$ws.Range("K132").Value = 16033.1661
$ws.Range("M132").Value = -13503.1661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 20.5
$ws.Range("I26").Value = 20.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 61.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 226.5
$ws.Range("N26").ClearContents()

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H68").Value = 1050
$ws.Range("I68").Value = 1050
$ws.Range("K68").Value = 3150
$ws.Range("M68").Value = -2339

$ws.Range("H70").Value = 8000
$ws.Range("J70").Value = 8000
$ws.Range("L70").Value = 24000
$ws.Range("N70").Value = -24630

$ws.Range("H71").Value = 1050
$ws.Range("I71").Value = 1050
$ws.Range("K71").Value = 9450
$ws.Range("M71").Value = -5394

$ws.Range("H73").Value = 8000
$ws.Range("J73").Value = 8000
$ws.Range("L73").Value = 24000
$ws.Range("N73").Value = -26184

$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 75000
$ws.Range("N81").Value = -77246

$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 225000
$ws.Range("N84").Value = -236232

$ws.Range("H131").Value = 906.5
$ws.Range("I131").Value = 610.6667
$ws.Range("J131").Value = 958.7059
$ws.Range("K131").Value = 1832.0001
$ws.Range("L131").Value = 2876.1177
$ws.Range("M131").Value = 3207.9999
$ws.Range("N131").Value = -12956.1177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15154
$ws.Range("J15").Value = 15154
$ws.Range("L15").Value = 15154
$ws.Range("N15").Value = -15730

$ws.Range("H81").Value = 15154
$ws.Range("J81").Value = 15154
$ws.Range("L81").Value = 15154
$ws.Range("N81").Value = -17150

$ws.Range("H84").Value = 15154
$ws.Range("J84").Value = 15154
$ws.Range("L84").Value = 45462
$ws.Range("N84").Value = -55446

$ws.Range("H102").Value = 1104.6
$ws.Range("I102").Value = 897.7857
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 897.7857
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = 724.2143
$ws.Range("N102").Value = -7244

$ws.Range("H123").Value = 39999
$ws.Range("J123").Value = 39999
$ws.Range("L123").Value = 39999
$ws.Range("N123").Value = -44899

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3538.2
$ws.Range("I7").Value = 2699
$ws.Range("K7").Value = 2699
$ws.Range("M7").Value = -2587

$ws.Range("H100").Value = 2279.2
$ws.Range("I100").Value = 2279.2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2279.2
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1738.2
$ws.Range("N100").ClearContents()

$ws.Range("H126").Value = 3538.2
$ws.Range("I126").Value = 2699
$ws.Range("K126").Value = 8097
$ws.Range("M126").Value = -5627

$ws.Range("H136").Value = 4866.5
$ws.Range("I136").Value = 4799.75
$ws.Range("K136").Value = 14399.25
$ws.Range("M136").Value = -11849.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
